# CHE_convelec_onshorewind.xlsx — "Preparation for transport" edit
#
# The commit adds a new parameter row ("capacity_to_activity") right
# after the existing "buildrate" row (row 9) in Sheet1, which pushes
# every row from 10 onward down by one. We reproduce that with a real
# row insert so every downstream row (and the trailing formatting-only
# row at the bottom of the sheet) shifts automatically, then populate
# the freshly inserted row with its data, and finally repair the
# AutoFilter range / _FilterDatabase defined name that Excel does not
# auto-extend on a plain row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row above row 10 — shifts rows 10:429 to 11:430,
#    and the worksheet dimension grows from L429 to L430 automatically.
$ws.Rows("10:10").Insert()

# 2) Populate the new row 10 with the "capacity_to_activity" entry.
$ws.Range("A10").Value = "CHE"
$ws.Range("B10").Value = "conv_elec_onshorewind"
$ws.Range("C10").Value = "capacity_to_activity"
$ws.Range("D10").Value = "constant"
$ws.Range("G10").Value = 0.001
$ws.Range("H10").Value = "GW/TWh"

# 3) Re-apply the AutoFilter over the now one-row-taller table
#    (A5:L849 -> A5:L850), and repair the hidden _FilterDatabase
#    defined name that backs it (not auto-updated by Insert()).
$ws.AutoFilterMode = $false
$ws.Range("A5:L850").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$5:`$L`$850"
    }
}

# 4) Match the author's final selection, sitting on the new value cell.
$ws.Range("G10").Select() | Out-Null
